$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-relevant header row (A1:D1) to new field names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case connector words (de/del/la/las/los/el/y) in state/municipality names
# and fix a stray double-capitalization ("MonteMorelos" -> "Montemorelos")
$ws.Range("B6").Value = "Pabellón De Arteaga"
$ws.Range("B7").Value = "Rincón De Romos"
$ws.Range("B24").Value = "Amatenango De La Frontera"
$ws.Range("B26").Value = "Bejucal De Ocampo"
$ws.Range("B33").Value = "Comitán De Domínguez"
$ws.Range("B51").Value = "Marqués De Comillas"
$ws.Range("B52").Value = "Mazapa De Madero"
$ws.Range("B57").Value = "Ocozocoautla De Espinosa"
$ws.Range("B63").Value = "San Cristóbal De Las Casas"
$ws.Range("B91").Value = "Hidalgo Del Parral"
$ws.Range("B95").Value = "San Francisco Del Oro"
$ws.Range("B106").Value = "San Juan De Sabinas"
$ws.Range("A114").Value = "Ciudad De México"
$ws.Range("B118").Value = "Cuajimalpa De Morelos"
$ws.Range("B141").Value = "San Juan Del Río"
$ws.Range("A149").Value = "Estado De México"
$ws.Range("B149").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B151").Value = "Almoloya De Alquisiras"
$ws.Range("B152").Value = "Almoloya De Juárez"
$ws.Range("B155").Value = "Atizapán De Zaragoza"
$ws.Range("B161").Value = "Coacalco De Berriozábal"
$ws.Range("B164").Value = "Ecatepec De Morelos"
$ws.Range("B167").Value = "Ixtapan Del Oro"
$ws.Range("B176").Value = "Naucalpan De Juárez"
$ws.Range("B180").Value = "San Felipe Del Progreso"
$ws.Range("B195").Value = "Tlalnepantla De Baz"
$ws.Range("B199").Value = "Valle De Bravo"
$ws.Range("B200").Value = "Valle De Chalco Solidaridad"
$ws.Range("B201").Value = "Villa De Allende"
$ws.Range("B211").Value = "Apaseo El Alto"
$ws.Range("B212").Value = "Apaseo El Grande"
$ws.Range("B219").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B223").Value = "Jaral Del Progreso"
$ws.Range("B232").Value = "San Diego De La Unión"
$ws.Range("B234").Value = "San Francisco Del Rincón"
$ws.Range("B236").Value = "San Luis De La Paz"
$ws.Range("B237").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B242").Value = "Valle De Santiago"
$ws.Range("B247").Value = "Acapulco De Juárez"
$ws.Range("B249").Value = "Ajuchitlán Del Progreso"
$ws.Range("B250").Value = "Alcozauca De Guerrero"
$ws.Range("B252").Value = "Atenango Del Río"
$ws.Range("B254").Value = "Atoyac De Álvarez"
$ws.Range("B255").Value = "Ayutla De Los Libres"
$ws.Range("B256").Value = "Chilapa De Álvarez"
$ws.Range("B257").Value = "Chilpancingo De Los Bravo"
$ws.Range("B258").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B262").Value = "Coyuca De Catalán"
$ws.Range("B266").Value = "Cutzamala De Pinzón"
$ws.Range("B269").Value = "Huitzuco De Los Figueroa"
$ws.Range("B270").Value = "Iguala De La Independencia"
$ws.Range("B272").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B273").Value = "Zihuatanejo De Azueta"
$ws.Range("B276").Value = "Mártir De Cuilapan"
$ws.Range("B288").Value = "Taxco De Alarcón"
$ws.Range("B290").Value = "Técpan De Galeana"
$ws.Range("B292").Value = "Tepecoacuilco De Trujano"
$ws.Range("B293").Value = "Tixtla De Guerrero"
$ws.Range("B297").Value = "Tlapa De Comonfort"
$ws.Range("B313").Value = "Cuautepec De Hinojosa"
$ws.Range("B316").Value = "Huejutla De Reyes"
$ws.Range("B321").Value = "Mixquiahuala De Juárez"
$ws.Range("B322").Value = "Molango De Escamilla"
$ws.Range("B324").Value = "Nopala De Villagrán"
$ws.Range("B325").Value = "Omitlán De Juárez"
$ws.Range("B326").Value = "Pachuca De Soto"
$ws.Range("B329").Value = "Progreso De Obregón"
$ws.Range("B336").Value = "Tenango De Doria"
$ws.Range("B337").Value = "Tepehuacán De Guerrero"
$ws.Range("B340").Value = "Tezontepec De Aldama"
$ws.Range("B344").Value = "Tula De Allende"
$ws.Range("B345").Value = "Tulancingo De Bravo"
$ws.Range("B347").Value = "Zacualtipán De Ángeles"
$ws.Range("B354").Value = "Atemajac De Brizuela"
$ws.Range("B356").Value = "Atotonilco El Alto"
$ws.Range("B357").Value = "Autlán De Navarro"
$ws.Range("B367").Value = "Encarnación De Díaz"
$ws.Range("B371").Value = "Ixtlahuacán Del Río"
$ws.Range("B378").Value = "Lagos De Moreno"
$ws.Range("B381").Value = "Ojuelos De Jalisco"
$ws.Range("B384").Value = "San Cristóbal De La Barranca"
$ws.Range("B385").Value = "San Juan De Los Lagos"
$ws.Range("B386").Value = "San Martín De Bolaños"
$ws.Range("B388").Value = "Santa María Del Oro"
$ws.Range("B390").Value = "Tamazula De Gordiano"
$ws.Range("B393").Value = "Tizapán El Alto"
$ws.Range("B394").Value = "Tlajomulco De Zúñiga"
$ws.Range("B399").Value = "Unión De San Antonio"
$ws.Range("B400").Value = "Valle De Guadalupe"
$ws.Range("B404").Value = "Yahualica De González Gallo"
$ws.Range("B405").Value = "Zacoalco De Torres"
$ws.Range("B407").Value = "Zapotlán El Grande"
$ws.Range("B481").Value = "Puente De Ixtla"
$ws.Range("B486").Value = "Tlaltizapán De Zapata"
$ws.Range("B493").Value = "Amatlán De Cañas"
$ws.Range("B496").Value = "Ixtlán Del Río"
$ws.Range("B500").Value = "Santa María Del Oro"
$ws.Range("B511").Value = "Mier Y Noriega"
$ws.Range("B512").Value = "Montemorelos"
$ws.Range("B514").Value = "San Nicolás De Los Garza"
$ws.Range("B516").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B518").Value = "Chalcatongo De Hidalgo"
$ws.Range("B520").Value = "Coicoyán De Las Flores"
$ws.Range("B522").Value = "El Barrio De La Soledad"
$ws.Range("B523").Value = "Fresnillo De Trujano"
$ws.Range("B524").Value = "Guadalupe De Ramírez"
$ws.Range("B526").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B527").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B528").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B531").Value = "Ixtlán De Juárez"
$ws.Range("B532").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B534").Value = "Mártires De Tacubaya"
$ws.Range("B537").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B539").Value = "Oaxaca De Juárez"
$ws.Range("B540").Value = "Ocotlán De Morelos"
$ws.Range("B541").Value = "Putla Villa De Guerrero"
$ws.Range("B548").Value = "San Baltazar Yatzachi El Bajo"
$ws.Range("B552").Value = "San Francisco Del Mar"
$ws.Range("B563").Value = "San José Del Progreso"
$ws.Range("B583").Value = "San Mateo Del Mar"
$ws.Range("B587").Value = "San Miguel Del Puerto"
$ws.Range("B593").Value = "San Pedro El Alto"
$ws.Range("B606").Value = "Santa Ana Del Valle"
$ws.Range("B650").Value = "Santo Domingo De Morelos"
$ws.Range("B660").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B662").Value = "Tataltepec De Valdés"
$ws.Range("B663").Value = "Teotitlán Del Valle"
$ws.Range("B664").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B665").Value = "Tlacolula De Matamoros"
$ws.Range("B666").Value = "Tlalixtac De Cabrera"
$ws.Range("B667").Value = "Totontepec Villa De Morelos"
$ws.Range("B670").Value = "Villa De Chilapa De Díaz"
$ws.Range("B671").Value = "Villa De Etla"
$ws.Range("B672").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B673").Value = "Villa De Zaachila"
$ws.Range("B674").Value = "Zapotitlán Del Río"
$ws.Range("B676").Value = "Zimatlán De Álvarez"
$ws.Range("B695").Value = "Chila De La Sal"
$ws.Range("B697").Value = "Cuetzalan Del Progreso"
$ws.Range("B706").Value = "Huehuetlán El Chico"
$ws.Range("B712").Value = "Izúcar De Matamoros"
$ws.Range("B719").Value = "Los Reyes De Juárez"
$ws.Range("B725").Value = "Palmar De Bravo"
$ws.Range("B738").Value = "San Salvador El Verde"
$ws.Range("B744").Value = "Tepanco De López"
$ws.Range("B745").Value = "Tepexi De Rodríguez"
$ws.Range("B747").Value = "Tetela De Ocampo"
$ws.Range("B757").Value = "Totoltepec De Guerrero"
$ws.Range("B770").Value = "Amealco De Bonfil"
$ws.Range("B772").Value = "Cadereyta De Montes"
$ws.Range("B776").Value = "Jalpan De Serra"
$ws.Range("B777").Value = "Landa De Matamoros"
$ws.Range("B779").Value = "Pinal De Amoles"
$ws.Range("B782").Value = "San Juan Del Río"
$ws.Range("B791").Value = "Axtla De Terrazas"
$ws.Range("B795").Value = "Ciudad Del Maíz"
$ws.Range("B803").Value = "Mexquitic De Carmona"
$ws.Range("B808").Value = "San Ciro De Acosta"
$ws.Range("B811").Value = "Santa María Del Río"
$ws.Range("B818").Value = "Tanquián De Escobedo"
$ws.Range("B820").Value = "Villa De Arriaga"
$ws.Range("B821").Value = "Villa De Guadalupe"
$ws.Range("B822").Value = "Villa De La Paz"
$ws.Range("B823").Value = "Villa De Ramos"
$ws.Range("B850").Value = "Jalpa De Méndez"
$ws.Range("B869").Value = "Soto La Marina"
$ws.Range("B883").Value = "San Pablo Del Monte"
$ws.Range("B885").Value = "Tetla De La Solidaridad"
$ws.Range("B897").Value = "Amatlán De Los Reyes"
$ws.Range("B901").Value = "Boca Del Río"
$ws.Range("B905").Value = "Cazones De Herrera"
$ws.Range("B913").Value = "Cosamaloapan De Carpio"
$ws.Range("B914").Value = "Cosautlán De Carvajal"
$ws.Range("B927").Value = "Hueyapan De Ocampo"
$ws.Range("B931").Value = "Ixhuatlán Del Café"
$ws.Range("B937").Value = "Juchique De Ferrer"
$ws.Range("B940").Value = "Lerdo De Tejada"
$ws.Range("B943").Value = "Martínez De La Torre"
$ws.Range("B952").Value = "Paso Del Macho"
$ws.Range("B955").Value = "Poza Rica De Hidalgo"
$ws.Range("B963").Value = "Soledad De Doblado"
$ws.Range("B967").Value = "Tatahuicapan De Juárez"
$ws.Range("B988").Value = "Vega De Alatorre"
$ws.Range("B1005").Value = "Cañitas De Felipe Pescador"
$ws.Range("B1006").Value = "Concepción Del Oro"
$ws.Range("B1015").Value = "Mezquital Del Oro"
$ws.Range("B1018").Value = "Moyahua De Estrada"
$ws.Range("B1019").Value = "Nochistlán De Mejía"
$ws.Range("B1027").Value = "Teúl De González Ortega"
$ws.Range("B1028").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B1031").Value = "Villa De Cos"

# Remove trailing footnote rows (1039:1043); the sheet now ends at row 1037
$ws.Rows("1039:1043").Delete() | Out-Null

Write-Host "edit complete"
